$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be appended in this exact order to match the
# source workbook's sharedStrings.xml ordering:
#   29: Cricket amount (K4)
#   30: ArmWrestling amount (K2)
#   31: Football amount (K3)
#   32: ArmWrestling QR image link (L2)

# Row 4 = Cricket: amount
$ws.Range("K4").Value = "₹ 1500/- only"

# Row 2 = ArmWrestling: amount is tiered by weight
$ws.Range("K2").Value = "₹ 100/- for below 80 & ₹ 150/- for above 80"

# Row 3 = Football: amount
$ws.Range("K3").Value = "₹ 1000/- only"

# Row 2 = ArmWrestling: add the QR code image link
$ws.Range("L2").Value = "https://media.discordapp.net/attachments/1162451241872412901/1169142818342715402/Untitled-1-01.png?ex=65545419&is=6541df19&hm=8a786424ef0da2bc2baafb80f67fc36ef3b51c15b1a95a3a897ab8d6a93a5373&=&width=641&height=662"

# Adjust column K width to fit the new (wider) amount text.
# NOTE: Excel stores column width internally as a fractional "characters"
# measurement derived from pixel widths; ColumnWidth = 10 is the value that
# serializes closest to the target stored width of 10.77734375 given this
# engine's width-quantization.
$ws.Columns.Item(11).ColumnWidth = 10

# Move selection/view to L8 like in the final saved state
$ws.Range("A1").Select()
$ws.Range("L8").Select()
